$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C numeric updates
$ws.Range("C2").Value = 17
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 15
$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 16
$ws.Range("C8").Value = 17
$ws.Range("C9").Value = 13
$ws.Range("C10").Value = 13
$ws.Range("C12").Value = 11
$ws.Range("C13").Value = 15
$ws.Range("C15").Value = 11
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 18
$ws.Range("C18").Value = 14

# Column B text updates
$ws.Range("B10").Value = "<hin>"
$ws.Range("B11").Value = "<mike>"
$ws.Range("B16").Value = "<zulu>"
$ws.Range("B18").Value = "<unifor>"
